# Updates the cryptos list (Coin/Link/Price/Volume) to reflect the latest
# pull from the data source. Prices stored in column D and volume deltas in
# column E are plain text in this sheet (not numbers/percentages), so we
# force the Text number format before writing and restore the default
# "Normal" style afterwards to avoid leaving a stray cell format behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.300.51"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.03%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.550.96"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.41%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "209.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.55%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.485"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.38%  "

# Row 7
$ws.Range("E7").Value = "  -0.05%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.67"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.28%  "

# Row 9
$ws.Range("E9").Value = "  -2.21%  "

# Row 10
$ws.Range("E10").Value = "  -1.35%  "

# Row 11
$ws.Range("E11").Value = "  +0.15%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.772.44"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.41%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.547.36"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.62%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.293.38"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.00%  "

# Row 15
$ws.Range("E15").Value = "  -1.47%  "

# Row 16
$ws.Range("E16").Value = "  -2.48%  "

# Row 17
$ws.Range("E17").Value = "  -2.89%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "227.06"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.62%  "

# Row 19
$ws.Range("E19").Value = "  -0.77%  "

# Row 20
$ws.Range("E20").Value = "  -2.73%  "

# Row 21
$ws.Range("E21").Value = "  -0.05%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.91"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.26%  "

# Row 23
$ws.Range("E23").Value = "  -3.20%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.06%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.71"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.16%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "14.75"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.67%  "

# Row 27
$ws.Range("E27").Value = "  -0.89%  "

# Row 28
$ws.Range("E28").Value = "  -0.12%  "

# Row 29
$ws.Range("E29").Value = "  -3.24%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0466"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.21%  "

# Row 31
$ws.Range("E31").Value = "  -4.51%  "

# Row 32
$ws.Range("E32").Value = "  -1.59%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.381.95"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.76%  "

# Row 35
$ws.Range("E35").Value = "  +0.97%  "

# Row 36
$ws.Range("E36").Value = "  -3.61%  "

# Row 37
$ws.Range("E37").Value = "  -1.30%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.59"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.14%  "

# Row 39
$ws.Range("E39").Value = "  -3.06%  "

# Row 40
$ws.Range("B40").Value = "ImmutableX"
$ws.Range("C40").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.510"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.71%  "

# Row 41
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.91"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.17%  "

# Row 42
$ws.Range("E42").Value = "  -0.09%  "

# Row 43
$ws.Range("E43").Value = "  -1.98%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0465"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.82%  "

# Row 46
$ws.Range("E46").Value = "  -2.03%  "

# Row 47
$ws.Range("E47").Value = "  -6.36%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.686.04"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.36%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "85.46"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.35%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "42.25"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.58%  "

# Row 51
$ws.Range("E51").Value = "  +0.55%  "

